$d = $word.ActiveDocument

$replacements = @(
    @("2025-01-28 Tuesday", "2025-01-29 Wednesday"),
    @("58×38=2204", "92×95=8740"),
    @("34×67=2278", "45×28=1260"),
    @("77×84=6468", "11×60=660"),
    @("99×39=3861", "98×79=7742"),
    @("12×94=1128", "45×38=1710"),
    @("80×59=4720", "57×58=3306"),
    @("78×13=1014", "40×38=1520"),
    @("34×87=2958", "37×89=3293"),
    @("47×98=4606", "54×49=2646"),
    @("89×97=8633", "79×96=7584"),
    @("38×63=2394", "80×66=5280"),
    @("63×14=882", "66×83=5478"),
    @("57×81=4617", "11×46=506"),
    @("65×91=5915", "59×59=3481"),
    @("12×18=216", "61×24=1464"),
    @("86×43=3698", "53×67=3551"),
    @("91×60=5460", "35×48=1680"),
    @("30×17=510", "25×93=2325"),
    @("96×31=2976", "44×94=4136"),
    @("16×51=816", "92×37=3404"),
    @("22×68=1496", "14×78=1092"),
    @("80×55=4400", "38×83=3154"),
    @("99×29=2871", "63×51=3213"),
    @("33×83=2739", "75×46=3450"),
    @("63×16=1008", "93×55=5115")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
